# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.488.06'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '1.570.86'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -1.54%  '
$ws.Range('D5').Formula = '="211.15"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  -1.53%  '
$ws.Range('D8').Formula = '="22.93"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +3.78%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D12').Value = '1.796.23'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '1.555.99'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '27.466.16'
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').Formula = '="62.35"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Formula = '="225.94"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +4.99%  '
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').Value = '0.0₃0705'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('E23').Value = '  +2.54%  '
$ws.Range('D24').Formula = '="1.95"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Formula = '="150.36"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -2.35%  '
$ws.Range('D26').Formula = '="15.16"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E28').Value = '  +1.37%  '
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').Value = '1.454.74'
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  -1.31%  '
$ws.Range('D38').Formula = '="0.0169"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').Formula = '="0.540"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('D40').Formula = '="0.813"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').Formula = '="2.36"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').Formula = '="0.992"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('D43').Formula = '="5.65"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -3.16%  '
$ws.Range('E44').Value = '  +6.73%  '
$ws.Range('E45').Value = '  -3.09%  '
$ws.Range('D46').Formula = '="63.88"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('D47').Value = '1.708.17'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').Formula = '="86.96"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').Formula = '="0.0947"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  -1.56%  '

$excel.CutCopyMode = 0
